$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.592.71"
$ws.Range("E2").Value = "  +3.81%  "
$ws.Range("D3").Value = "3.091.84"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'218.87"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").Value = "'617.18"
$ws.Range("E6").Value = "  -2.94%  "
$ws.Range("D7").Value = "'0.375"
$ws.Range("E7").Value = "  -3.18%  "
$ws.Range("D8").Value = "'0.921"
$ws.Range("E8").Value = "  +11.64%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "3.090.99"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").Value = "'0.684"
$ws.Range("E11").Value = "  +17.36%  "
$ws.Range("D12").Value = "'0.190"
$ws.Range("E12").Value = "  +6.13%  "
$ws.Range("D13").Value = "'0.0000254"
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").Value = "91.175.42"
$ws.Range("E14").Value = "  +3.56%  "
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "'33.04"
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("D17").Value = "3.652.97"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "3.134.18"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'3.50"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("D20").Value = "'0.0000220"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'13.77"
$ws.Range("E21").Value = "  +3.29%  "
$ws.Range("D22").Value = "'435.06"
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").Value = "'8.48"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +4.76%  "
$ws.Range("D25").Value = "'5.59"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").Value = "'84.34"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "'11.80"
$ws.Range("E27").Value = "  +2.68%  "
$ws.Range("D28").Value = "3.257.49"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'0.168"
$ws.Range("E30").Value = "  +7.24%  "
$ws.Range("D31").Value = "'1.01"
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("D32").Value = "'8.76"
$ws.Range("E32").Value = "  +6.89%  "
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").Value = "'518.95"
$ws.Range("E34").Value = "  +2.91%  "
$ws.Range("D35").Value = "'7.05"
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -7.36%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "'23.04"
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("D40").Value = "'22.33"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "'0.370"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.88"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("D46").Value = "'0.0726"
$ws.Range("E46").Value = "  +7.44%  "
$ws.Range("D47").Value = "'43.95"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'141.45"
$ws.Range("E48").Value = "  -3.26%  "
$ws.Range("D49").Value = "'0.000263"
$ws.Range("E49").Value = "  +12.13%  "
$ws.Range("E50").Value = "  +5.90%  "
$ws.Range("D51").Value = "'164.16"
$ws.Range("E51").Value = "  +1.00%  "
